$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J4").Value = 1869
$ws.Range("J7").Value = 29345
$ws.Range("L2").Value = 3896
$ws.Range("L3").Value = 4071
$ws.Range("L4").Value = 1002
$ws.Range("L5").Value = 237
$ws.Range("L6").Value = 3526
$ws.Range("L7").Value = 12732

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J63").Value = 221
$ws.Range("J101").Value = 29345
$ws.Range("K63").Value = 166
$ws.Range("K65").Value = 637
$ws.Range("L6").Value = 104
$ws.Range("L8").Value = 831
$ws.Range("L10").Value = 80
$ws.Range("L11").Value = 210
$ws.Range("L15").Value = 93
$ws.Range("L19").Value = 368
$ws.Range("L20").Value = 317
$ws.Range("L23").Value = 135
$ws.Range("L25").Value = 70
$ws.Range("L29").Value = 703
$ws.Range("L33").Value = 596
$ws.Range("L36").Value = 164
$ws.Range("L37").Value = 467
$ws.Range("L42").Value = 408
$ws.Range("L43").Value = 97
$ws.Range("L44").Value = 91
$ws.Range("L47").Value = 93
$ws.Range("L51").Value = 157
$ws.Range("L52").Value = 256
$ws.Range("L53").Value = 142
$ws.Range("L54").Value = 264
$ws.Range("L55").Value = 119
$ws.Range("L63").Value = 42
$ws.Range("L64").Value = 85
$ws.Range("L65").Value = 238
$ws.Range("L67").Value = 451
$ws.Range("L71").Value = 36
$ws.Range("L76").Value = 193
$ws.Range("L78").Value = 160
$ws.Range("L79").Value = 333
$ws.Range("L84").Value = 123
$ws.Range("L85").Value = 669
$ws.Range("L88").Value = 141
$ws.Range("L89").Value = 178
$ws.Range("L93").Value = 69
$ws.Range("L96").Value = 134
$ws.Range("L101").Value = 12732

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L4").Value = 16
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 178

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 201
$ws.Range("L4").Value = 45
$ws.Range("L7").Value = 669

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 88
$ws.Range("L3").Value = 78
$ws.Range("L7").Value = 256

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 142

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 244
$ws.Range("L3").Value = 277
$ws.Range("L6").Value = 226
$ws.Range("L7").Value = 831

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 165
$ws.Range("L6").Value = 190
$ws.Range("L7").Value = 596

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 155
$ws.Range("L4").Value = 28
$ws.Range("L6").Value = 128
$ws.Range("L7").Value = 467

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K4").Value = 29
$ws.Range("K7").Value = 637
$ws.Range("L3").Value = 70
$ws.Range("L7").Value = 238

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L6").Value = 105
$ws.Range("L7").Value = 451

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 46
$ws.Range("L4").Value = 5
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 125
$ws.Range("L7").Value = 264

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 219
$ws.Range("L3").Value = 263
$ws.Range("L7").Value = 703

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 130
$ws.Range("L3").Value = 112
$ws.Range("L6").Value = 107
$ws.Range("L7").Value = 368

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 39
$ws.Range("L6").Value = 88
$ws.Range("L7").Value = 193

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 132
$ws.Range("L7").Value = 408

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 80

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L4").Value = 17
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 39
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 119

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 34
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 121
$ws.Range("L6").Value = 69
$ws.Range("L7").Value = 333

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L6").Value = 82
$ws.Range("L7").Value = 317

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 24
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L2").Value = 24
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 33
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L5").Value = 3
$ws.Range("L7").Value = 141

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 48
$ws.Range("L7").Value = 157

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L2").Value = 16
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 36
